$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (single-dot decimal-looking strings),
# so they remain text cells exactly like the surrounding inlineStr cells.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"

# Apply the updated Price (D) / Volume(1h) (E) values from the source refresh.
$ws.Range('D2').Value = '69.668.11'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = '2.504.44'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '574.32'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').Value = '166.53'
$ws.Range('E6').Value = '  -1.00%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('D9').Value = '2.502.76'
$ws.Range('E9').Value = '  -0.65%  '
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('D12').Value = '0.358'
$ws.Range('E12').Value = '  +2.70%  '
$ws.Range('D13').Value = '4.94'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').Value = '2.959.39'
$ws.Range('D15').Value = '69.519.27'
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').Value = '0.0000176'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').Value = '24.72'
$ws.Range('E17').Value = '  -1.77%  '
$ws.Range('D18').Value = '2.501.89'
$ws.Range('D19').Value = '11.21'
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('E20').Value = '  -3.60%  '
$ws.Range('D21').Value = '349.15'
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('D22').Value = '3.90'
$ws.Range('E22').Value = '  -1.03%  '
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('D25').Value = '70.80'
$ws.Range('E25').Value = '  +2.15%  '
$ws.Range('E26').Value = '  -2.05%  '
$ws.Range('D27').Value = '8.74'
$ws.Range('E27').Value = '  -3.35%  '
$ws.Range('D28').Value = '2.630.26'
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '0.0₃0890'
$ws.Range('E30').Value = '  -2.18%  '
$ws.Range('E31').Value = '  -0.59%  '
$ws.Range('D32').Value = '458.15'
$ws.Range('E32').Value = '  -1.79%  '
$ws.Range('E33').Value = '  -5.95%  '
$ws.Range('E34').Value = '  -1.61%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').Value = '157.29'
$ws.Range('E36').Value = '  +2.66%  '
$ws.Range('E37').Value = '  -3.42%  '
$ws.Range('D39').Value = '18.36'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('E41').Value = '  -1.29%  '
$ws.Range('D42').Value = '4.69'
$ws.Range('E42').Value = '  -2.15%  '
$ws.Range('D43').Value = '1.61'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').Value = '38.15'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('E45').Value = '  -5.27%  '
$ws.Range('E46').Value = '  -7.89%  '
$ws.Range('D47').Value = '141.34'
$ws.Range('E47').Value = '  -1.45%  '
$ws.Range('E48').Value = '  -0.61%  '
$ws.Range('E49').Value = '  -2.64%  '
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('E51').Value = '  -0.43%  '
